$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert 4 new columns before the old column H ---------------------
# This shifts old H,I,J (csv/xls filename columns) to L,M,N
$ws.Range("H1:K1").EntireColumn.Insert()

# --- 2) Fill in the headers for the newly inserted columns ----------------
$ws.Range("H1").Value = "start_angle_standard"
$ws.Range("I1").Value = "end_angle_standard"
$ws.Range("J1").Value = "start_angle_relax"
$ws.Range("K1").Value = "end_angle_relax"

# Match formatting of the surrounding header cells (centered header style)
$ws.Range("G1").Copy()
$ws.Range("H1:K1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3) Clear out row 3 entirely, keep only a formatted empty cell at N3 --
$ws.Range("A3:N3").ClearContents()
$ws.Range("N3").NumberFormat = "@"

# Make sure the data-file columns keep their original text format
$ws.Range("N2").NumberFormat = "@"

# --- 4) Column widths (converted to Excel's character-width model) --------
$ws.Range("A1").ColumnWidth = 5.428571428571429
$ws.Range("B1").ColumnWidth = 12.142857142857142
$ws.Range("C1").ColumnWidth = 13.428571428571429
$ws.Range("D1").ColumnWidth = 11.0
$ws.Range("E1").ColumnWidth = 14.142857142857142
$ws.Range("F1").ColumnWidth = 22.428571428571427
$ws.Range("G1").ColumnWidth = 19.714285714285715
$ws.Range("H1").ColumnWidth = 21.571428571428573
$ws.Range("I1").ColumnWidth = 20.428571428571427
$ws.Range("J1").ColumnWidth = 19.285714285714285
$ws.Range("K1").ColumnWidth = 22.714285714285715
$ws.Range("L1").ColumnWidth = 27.857142857142858
$ws.Range("M1").ColumnWidth = 27.571428571428573
$ws.Range("N1").ColumnWidth = 18.285714285714285

# --- 5) Selection, matching the final state recorded in the workbook ------
[void]$ws.Range("M19").Select()

[void]$wb.Save()
